$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A510").Value = "Efx"
$ws.Range("B510").Value = "sftflf"
$ws.Range("C510").Value = 2
$ws.Range("D510").Value = "longitude latitude typefis"
$ws.Range("E510").Value = "Floating Ice Shelf  Area Fraction"
$ws.Range("F510").Value = "%"
$ws.Range("G510").Value = "web"
$ws.Range("H510").Value = "Not available in LIM"
$ws.Range("I510").Value = "David Docquier, Thomas"
$ws.Range("J510").Value = "Fraction of grid cell covered by floating ice shelf, the component of the ice sheet that is flowing over sea water"
$ws.Range("K510").Value = "CMIP,ISMIP6"

$ws.Range("A512").Value = "SImon"
$ws.Range("B512").Value = "simpconc"
$ws.Range("C512").Value = 3
$ws.Range("D512").Value = "longitude latitude time typemp"
$ws.Range("E512").Value = "Percentage Cover of Sea-Ice by Meltpond"
$ws.Range("F512").Value = "%"
$ws.Range("G512").Value = "web"
$ws.Range("H512").Value = "Not available in LIM in EC-Earth3's CMIP6 version. In a newer version: simpconc = iceamp / siconc  according to David, and probably not  simpconc =  ( iceamp (= melt-pond fraction per grid-cell area, no unit) / grid-cell area ) * 100 [in %]"
$ws.Range("I512").Value = "David Docquier, Thomas"
$ws.Range("J512").Value = "Percentage of sea ice, by area, which is covered by melt ponds, giving equal weight to every square metre of sea ice ."
$ws.Range("K512").Value = "C4MIP,CMIP,FAFMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,RFMIP,SIMIP"

$ws.Range("A513").Value = "SImon"
$ws.Range("B513").Value = "sirdgconc"
$ws.Range("C513").Value = 3
$ws.Range("D513").Value = "longitude latitude time typesirdg"
$ws.Range("E513").Value = "Percentage Cover of Sea-Ice by Ridging"
$ws.Range("F513").Value = 1
$ws.Range("G513").Value = "web"
$ws.Range("H513").Value = "Not available in LIM"
$ws.Range("I513").Value = "David Docquier, Thomas"
$ws.Range("J513").Value = "Fraction of sea ice, by area, which is covered by sea ice ridges, giving equal weight to every square metre of sea ice ."
$ws.Range("K513").Value = "C4MIP,CMIP,FAFMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,RFMIP,SIMIP"

$ws.Range("A514").Value = "SImon"
$ws.Range("B514").Value = "sipr"
$ws.Range("C514").Value = 2
$ws.Range("D514").Value = "longitude latitude time"
$ws.Range("E514").Value = "Rainfall rate over sea ice"
$ws.Range("F514").Value = "kg m-2 s-1"
$ws.Range("G514").Value = "web"
$ws.Range("H514").Value = "Not available in LIM"
$ws.Range("I514").Value = "David Docquier, Thomas"
$ws.Range("J514").Value = "mass of liquid precipitation falling onto sea ice divided by grid-cell area"
$ws.Range("K514").Value = "C4MIP,CFMIP,CMIP,FAFMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,PMIP,RFMIP,SIMIP,VIACSAB"

$ws.Range("A516").Value = "LImon"
$ws.Range("B516").Value = "sftflf"
$ws.Range("C516").Value = 2
$ws.Range("D516").Value = "longitude latitude time typefis"
$ws.Range("E516").Value = "Floating Ice Shelf  Area Fraction"
$ws.Range("F516").Value = "%"
$ws.Range("G516").Value = "web"
$ws.Range("H516").Value = "Not available in LIM"
$ws.Range("I516").Value = "David Docquier, Thomas"
$ws.Range("J516").Value = "Fraction of grid cell covered by floating ice shelf, the component of the ice sheet that is flowing over sea water"
$ws.Range("K516").Value = "CMIP,ISMIP6"

$ws.Range("A520").Value = "Omon"
$ws.Range("B520").Value = "dissi14cabioos"
$ws.Range("C520").Value = 2
$ws.Range("D520").Value = "longitude latitude time"
$ws.Range("E520").Value = "mole_concentration_of_dissolved_inorganic_carbon14_abiotic_analogue_in_sea_water"
$ws.Range("F520").Value = "mol m-3"
$ws.Range("G520").Value = "web"
$ws.Range("H520").Value = "Not available"
$ws.Range("I520").Value = "Raffaele Bernardello"
$ws.Range("J520").Value = "Abiotic Dissolved inorganic 14carbon (CO3+HCO3+H2CO3) concentration"
$ws.Range("K520").Value = "AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP"

$ws.Range("A521").Value = "Omon"
$ws.Range("B521").Value = "vsfcorr"
$ws.Range("C521").Value = 2
$ws.Range("D521").Value = "longitude latitude time"
$ws.Range("E521").Value = "Virtual Salt Flux Correction"
$ws.Range("F521").Value = "kg m-2 s-1"
$ws.Range("G521").Value = "web"
$ws.Range("H521").Value = "Not available"
$ws.Range("I521").Value = "Raffaele Bernardello"
$ws.Range("J521").Value = "It is set to zero in models which receive a real water flux."
$ws.Range("K521").Value = "AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB"

$ws.Range("A522").Value = "Omon"
$ws.Range("B522").Value = "intppdiaz"
$ws.Range("C522").Value = 3
$ws.Range("D522").Value = "longitude latitude time"
$ws.Range("E522").Value = "Net Primary Mole Productivity of Carbon by Diazotrophs"
$ws.Range("F522").Value = "mol m-2 s-1"
$ws.Range("G522").Value = "web"
$ws.Range("H522").Value = "Not available: INTNFIX is the production by Diazotrophs but they do not contribute to carbon so I think this one is missing. "
$ws.Range("I522").Value = "Raffaele Bernardello, Thomas"
$ws.Range("J522").Value = "''Production of carbon' means the production of biomass expressed as the mass of carbon which it contains. Net primary production is the excess of gross primary production (rate of synthesis of biomass from inorganic precursors) by autotrophs ('producers'), for example, photosynthesis in plants or phytoplankton, over the rate at which the autotrophs themselves respire some of this biomass. 'Productivity' means production per unit area. In ocean modelling, diazotrophs are phytoplankton of the phylum cyanobacteria distinct from other phytoplankton groups in their ability to fix nitrogen gas in addition to nitrate and ammonium. Phytoplankton are autotrophic prokaryotic or eukaryotic algae that live near the water surface where there is sufficient light to support photosynthesis. The phrase 'expressed_as' is used in the construction A_expressed_as_B, where B is a chemical constituent of A."
$ws.Range("K522").Value = "AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB"

$ws.Range("A523").Value = "Omon"
$ws.Range("B523").Value = "intpppico"
$ws.Range("C523").Value = 3
$ws.Range("D523").Value = "longitude latitude time"
$ws.Range("E523").Value = "Net Primary Mole Productivity of Carbon by Picophytoplankton"
$ws.Range("F523").Value = "mol m-2 s-1"
$ws.Range("G523").Value = "web"
$ws.Range("H523").Value = "Not available"
$ws.Range("I523").Value = "Raffaele Bernardello"
$ws.Range("J523").Value = "''Production of carbon' means the production of biomass expressed as the mass of carbon which it contains. Net primary production is the excess of gross primary production (rate of synthesis of biomass from inorganic precursors) by autotrophs ('producers'), for example, photosynthesis in plants or phytoplankton, over the rate at which the autotrophs themselves respire some of this biomass. 'Productivity' means production per unit area. Picophytoplankton are phytoplankton of less than 2 micrometers in size. Phytoplankton are autotrophic prokaryotic or eukaryotic algae that live near the water surface where there is sufficient light to support photosynthesis. The phrase 'expressed_as' is used in the construction A_expressed_as_B, where B is a chemical constituent of A. It means that the quantity indicated by the standard name is calculated solely with respect to the B contained in A, neglecting all other chemical constituents of A."
$ws.Range("K523").Value = "AerChemMIP,C4MIP,CMIP,GMMIP,GeoMIP,HighResMIP,LS3MIP,OMIP,VIACSAB"

$ws.Range("A525").Value = "Oyr"
$ws.Range("B525").Value = "dissi14cabio"
$ws.Range("C525").Value = 1
$ws.Range("D525").Value = "longitude latitude olevel time"
$ws.Range("E525").Value = "Abiotic Dissolved Inorganic 14Carbon Concentration"
$ws.Range("F525").Value = "mol m-3"
$ws.Range("G525").Value = "web"
$ws.Range("H525").Value = "Not available"
$ws.Range("I525").Value = "Raffaele Bernardello"
$ws.Range("J525").Value = "Abiotic Dissolved inorganic 14carbon (CO3+HCO3+H2CO3) concentration"
$ws.Range("K525").Value = "AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP"

$ws.Range("A526").Value = "Oyr"
$ws.Range("B526").Value = "bacc"
$ws.Range("C526").Value = 3
$ws.Range("D526").Value = "longitude latitude olevel time"
$ws.Range("E526").Value = "Bacterial Carbon Concentration"
$ws.Range("F526").Value = "mol m-3"
$ws.Range("G526").Value = "web"
$ws.Range("H526").Value = "Not available"
$ws.Range("I526").Value = "Raffaele Bernardello"
$ws.Range("J526").Value = "Sum of bacterial carbon component concentrations"
$ws.Range("K526").Value = "AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP,PMIP,VIACSAB"

$ws.Range("A527").Value = "Oyr"
$ws.Range("B527").Value = "arag"
$ws.Range("C527").Value = 2
$ws.Range("D527").Value = "longitude latitude olevel time"
$ws.Range("E527").Value = "Aragonite Concentration"
$ws.Range("F527").Value = "mol m-3"
$ws.Range("G527").Value = "web"
$ws.Range("H527").Value = "Not available"
$ws.Range("I527").Value = "Raffaele Bernardello"
$ws.Range("J527").Value = "Sum of particulate aragonite components (e.g. Phytoplankton, Detrital, etc.)"
$ws.Range("K527").Value = "AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP,PMIP,VIACSAB"

$ws.Range("A528").Value = "Oyr"
$ws.Range("B528").Value = "phydiaz"
$ws.Range("C528").Value = 3
$ws.Range("D528").Value = "longitude latitude olevel time"
$ws.Range("E528").Value = "Mole Concentration of Diazotrophs expressed as Carbon in sea water"
$ws.Range("F528").Value = "mol m-3"
$ws.Range("G528").Value = "web"
$ws.Range("H528").Value = "Not available"
$ws.Range("I528").Value = "Raffaele Bernardello"
$ws.Range("J528").Value = "carbon concentration from the diazotrophic phytoplankton component alone"
$ws.Range("K528").Value = "AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP,PMIP,VIACSAB"

$ws.Range("A529").Value = "Oyr"
$ws.Range("B529").Value = "phycalc"
$ws.Range("C529").Value = 3
$ws.Range("D529").Value = "longitude latitude olevel time"
$ws.Range("E529").Value = "Mole Concentration of Calcareous Phytoplankton expressed as Carbon in sea water"
$ws.Range("F529").Value = "mol m-3"
$ws.Range("G529").Value = "web"
$ws.Range("H529").Value = "Not available"
$ws.Range("I529").Value = "Raffaele Bernardello"
$ws.Range("J529").Value = "carbon concentration from calcareous (calcite-producing) phytoplankton component alone"
$ws.Range("K529").Value = "AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP,PMIP,VIACSAB"

$ws.Range("A530").Value = "Oyr"
$ws.Range("B530").Value = "zmisc"
$ws.Range("C530").Value = 3
$ws.Range("D530").Value = "longitude latitude olevel time"
$ws.Range("E530").Value = "Mole Concentration of Other Zooplankton expressed as Carbon in sea water"
$ws.Range("F530").Value = "mol m-3"
$ws.Range("G530").Value = "web"
$ws.Range("H530").Value = "Not available"
$ws.Range("I530").Value = "Raffaele Bernardello"
$ws.Range("J530").Value = "carbon from additional zooplankton component concentrations alone (e.g. Micro, meso).  Since the models all have different numbers of components, this variable has been included to provide a check for intercomparison between models since some phytoplankton groups are supersets."
$ws.Range("K530").Value = "AerChemMIP,C4MIP,CMIP,GeoMIP,LUMIP,OMIP,PMIP,VIACSAB"

$ws.Range("A531").Value = "Oyr"
$ws.Range("B531").Value = "co3satarag"
$ws.Range("C531").Value = 2
$ws.Range("D531").Value = "longitude latitude olevel time"
$ws.Range("E531").Value = "Mole Concentration of Carbonate Ion in Equilibrium with Pure Aragonite in sea water"
$ws.Range("F531").Value = "mol m-3"
$ws.Range("G531").Value = "web"
$ws.Range("H531").Value = "Not available"
$ws.Range("I531").Value = "Raffaele Bernardello"
$ws.Range("J531").Value = "Mole concentration means number of moles per unit volume, also called 'molarity', and is used in the construction 'mole_concentration_of_X_in_Y', where X is a material constituent of Y. A chemical or biological species denoted by X may be described by a single term such as 'nitrogen' or a phrase such as 'nox_expressed_as_nitrogen'. The phrase 'expressed_as' is used in the construction A_expressed_as_B, where B is a chemical constituent of A. It means that the quantity indicated by the standard name is calculated solely with respect to the B contained in A, neglecting all other chemical constituents of A. The chemical formula of the carbonate anion is CO3 with an electrical charge of minus two. Aragonite is a mineral that is a polymorph of calcium carbonate. The chemical formula of aragonite is CaCO3. At a given salinity, the thermodynamic equilibrium is that between dissolved carbonate ion and solid aragonite. Standard names also exist for calcite, another polymorph of calcium carbonate."
$ws.Range("K531").Value = "AerChemMIP,CMIP,GeoMIP,LUMIP,OMIP,VIACSAB"

$ws.Range("A533").Value = "Ofx"
$ws.Range("B533").Value = "ugrid"
$ws.Range("C533").Value = 1
$ws.Range("D533").Value = "longitude latitude"
$ws.Range("E533").Value = "UGRID Grid Information"
$ws.Range("G533").Value = "web"
$ws.Range("H533").Value = "Not required because the NEMO grid is curvilinear"
$ws.Range("I533").Value = "Thomas Reerink"
$ws.Range("J533").Value = "Ony required for models with unstructured grids: this label should be used for a file containing information about the grid structure, following the UGRID convention."
$ws.Range("K533").Value = "CMIP"

$ws.Range("A536:Y539").Select()